$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$ch = $ws.Range("B8").Characters()
Write-Host ("Characters obj: " + $ch)
$ch.Text = "2025-07-29T07:08:53+00:00"
$v = $ws.Range("B8").Value()
Write-Host "B8: [$v]"
